$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header C3: "RMSE (eV)" -> "Testing RMSE (eV)"
$ws.Range("C3").Value = "Testing RMSE (eV)"

# Row 8 (A8:G8) already uses the bold/red font (style index 3/4 in the old file);
# update that font to bold + automatic/theme text color instead of red.
# Do this before touching alignment so A8 still shares the same font shape as F8:G8.
$ws.Range("A8:G8").Font.ThemeColor = 1
$ws.Range("A8:G8").Font.Bold = $true

# A4:A7 were numbered 1..4 -> become letters a..d, right-aligned (non-bold)
$ws.Range("A4").Value = "a"
$ws.Range("A5").Value = "b"
$ws.Range("A6").Value = "c"
$ws.Range("A7").Value = "d"
$ws.Range("A4:A7").HorizontalAlignment = -4152

# A8 was numbered 5 -> becomes letter e, also right-aligned
$ws.Range("A8").Value = "e"
$ws.Range("A8").HorizontalAlignment = -4152

Write-Output "done"
